# feat: add 2022-Q3 data
#
# The workbook has two sheets: "总计" (summary) and "2022-Q2" (fund
# holdings for that quarter). We need to:
#   1. Introduce a new "2022-Q3" sheet with fund-holdings data, inserted
#      right after "总计" (so "2022-Q2" ends up third).
#   2. Add a new top row to "总计" summarizing the 2022-Q3 data, pushing
#      the existing 2022-Q2 summary row down.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)
$wsQ2    = $wb.Worksheets.Item(2)

# --- Step 1: duplicate "2022-Q2" so its data survives untouched under its
# own sheet, then repurpose the original sheet (which keeps its rId/sheetId
# and thus its tab position right after "总计") for the new "2022-Q3" data.
$wsQ2.Copy($null, $wsQ2)
$wsQ2Copy = $wb.Worksheets.Item(3)

$wsQ2.Name = "2022-Q3"
$wsQ2Copy.Name = "2022-Q2"

$wsQ3 = $wsQ2

# --- Step 2: clear out the old Q2 fund-holdings content from $wsQ3 and
# write the new Q3 fund-holdings content in its place.
$wsQ3.Cells.Clear()

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Copy the header/"index" cell format (bold, centered, bordered) from the
# summary sheet's A2 onto the new sheet's header row and index column - the
# same style used throughout this workbook for header cells.
$wsTotal.Range("A2").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A4").PasteSpecial(-4122)

$wsQ3.Range("A2").Value = 0
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("A4").Value = 2

$wsQ3.Range("H2").Value = 10
$wsQ3.Range("H3").Value = 1
$wsQ3.Range("H4").Value = 1

$wsQ3.Range("C2").Value = "光大保德信一带一路战略主题混合"
$wsQ3.Range("C3").Value = "东方多策略灵活配置混合C"
$wsQ3.Range("C4").Value = "东方多策略灵活配置混合A"

# Columns B, D, E, F, G hold numeric-looking text ("001463", "1.57", ...)
# that must stay text (leading zeros, fixed decimal formatting) rather than
# become real numbers, so format them as Text before writing the values.
$wsQ3.Range("B2:B4").NumberFormat = "@"
$wsQ3.Range("D2:G4").NumberFormat = "@"

$wsQ3.Range("B2").Value = "001463"
$wsQ3.Range("D2").Value = "1.57"
$wsQ3.Range("E2").Value = "87.51"
$wsQ3.Range("F2").Value = "4.81"
$wsQ3.Range("G2").Value = "0.0755"

$wsQ3.Range("B3").Value = "002068"
$wsQ3.Range("D3").Value = "0.26"
$wsQ3.Range("E3").Value = "55.14"
$wsQ3.Range("F3").Value = "2.97"
$wsQ3.Range("G3").Value = "0.0077"

$wsQ3.Range("B4").Value = "400023"
$wsQ3.Range("D4").Value = "0.03"
$wsQ3.Range("E4").Value = "55.14"
$wsQ3.Range("F4").Value = "2.97"
$wsQ3.Range("G4").Value = "0.0009"

# --- Step 3: insert the 2022-Q3 summary row into "总计", above the
# existing 2022-Q2 row (copying its format so the new & shifted rows both
# keep the bold/bordered index-column style).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.08
